# Updates the "Price" (column D) and "Volume(1h)" (column E) figures in the
# cryptos list, matching the latest scrape pulled in by the scheduled
# GitHub Actions run.
#
# Notes:
#  - Column D/E cells are stored as text (inlineStr) in the source sheet,
#    including values that look numeric (e.g. "1.00", "69.00", "4.03").
#    Assigning such a string straight to .Value lets Excel's COM layer
#    auto-coerce it to a real number (dropping the trailing/insignificant
#    zeros and any formatting), so for any new Price value that Excel would
#    interpret as a plain number we prefix it with a leading apostrophe
#    (the standard Excel "force text" marker) before assigning .Value.
#    Values that already contain two dots (e.g. "67.189.01") are never
#    number-like, so they are assigned as-is.
#  - Percent cells in column E are always left as literal text (they carry
#    leading/trailing padding spaces), so no text-forcing is required there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.189.01"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").Value = "2.482.57"

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'585.27"

$ws.Range("D6").Value = "'172.77"
$ws.Range("E6").Value = "  +3.43%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.33%  "

$ws.Range("D9").Value = "2.483.18"
$ws.Range("E9").Value = "  +0.56%  "

$ws.Range("E10").Value = "  +3.08%  "

$ws.Range("E11").Value = "  +1.15%  "

$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("D13").Value = "'0.333"
$ws.Range("E13").Value = "  +0.16%  "

$ws.Range("D15").Value = "'25.56"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("D16").Value = "67.065.57"
$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("D17").Value = "'0.0000170"
$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").Value = "2.430.91"
$ws.Range("E18").Value = "  -1.40%  "

$ws.Range("D19").Value = "'7.58"
$ws.Range("E19").Value = "  -0.99%  "

$ws.Range("D20").Value = "'10.97"
$ws.Range("E20").Value = "  -3.41%  "

$ws.Range("D21").Value = "'350.47"
$ws.Range("E21").Value = "  -1.66%  "

$ws.Range("D22").Value = "'4.03"
$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").Value = "'69.00"
$ws.Range("E24").Value = "  -0.82%  "

$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("E26").Value = "  +2.93%  "

$ws.Range("E27").Value = "  +2.11%  "

$ws.Range("D28").Value = "2.608.61"
$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.40%  "

$ws.Range("E30").Value = "  +1.49%  "

$ws.Range("D31").Value = "'506.87"
$ws.Range("E31").Value = "  -1.12%  "

$ws.Range("E32").Value = "  -1.25%  "

$ws.Range("E33").Value = "  +1.54%  "

$ws.Range("E34").Value = "  -0.82%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").Value = "'162.36"
$ws.Range("E36").Value = "  +2.28%  "

$ws.Range("D37").Value = "'0.119"
$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("D39").Value = "'18.16"
$ws.Range("E39").Value = "  -1.65%  "

$ws.Range("E40").Value = "  -0.40%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("E42").Value = "  +1.03%  "

$ws.Range("E43").Value = "  +1.35%  "

$ws.Range("E44").Value = "  +0.94%  "

$ws.Range("E45").Value = "  +3.14%  "

$ws.Range("D46").Value = "'143.39"
$ws.Range("E46").Value = "  +1.33%  "

$ws.Range("E47").Value = "  +3.54%  "

$ws.Range("D50").Value = "'0.0737"
$ws.Range("E50").Value = "  +0.74%  "

$ws.Range("E51").Value = "  -0.74%  "
